$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.132.15"
$ws.Range("E2").Value = "  +1.11%  "
$ws.Range("D3").Value = "1.891.25"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.21"
$ws.Range("E5").Value = "  +1.35%  "
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5160"
$ws.Range("E7").Value = "  +2.35%  "
$ws.Range("E8").Value = "  +1.85%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07213"
$ws.Range("E9").Value = "  +0.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9047"
$ws.Range("E10").Value = "  +1.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.05"
$ws.Range("E11").Value = "  +1.81%  "
$ws.Range("E12").Value = "  +1.48%  "
$ws.Range("D13").Value = "1.887.20"
$ws.Range("E13").Value = "  +1.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.96"
$ws.Range("E14").Value = "  +2.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.276"
$ws.Range("E15").Value = "  +0.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9999"
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.36"
$ws.Range("E18").Value = "  +2.18%  "
$ws.Range("D20").Value = "27.169.77"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.055"
$ws.Range("E21").Value = "  +0.57%  "
$ws.Range("D22").Value = "2.125.78"
$ws.Range("E22").Value = "  +1.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.60"
$ws.Range("E23").Value = "  +2.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.433"
$ws.Range("E24").Value = "  -0.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.82"
$ws.Range("E25").Value = "  -0.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.793"
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.04"
$ws.Range("E27").Value = "  +1.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.152"
$ws.Range("E28").Value = "  +4.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.74"
$ws.Range("E29").Value = "  +1.75%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.991"
$ws.Range("E30").Value = "  +7.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.816"
$ws.Range("E31").Value = "  +3.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09214"
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05065"
$ws.Range("E33").Value = "  -0.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.202"
$ws.Range("E34").Value = "  +4.85%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7594"
$ws.Range("E35").Value = "  +2.42%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.021"
$ws.Range("E36").Value = "  +1.58%  "
$ws.Range("E37").Value = "  +1.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.566"
$ws.Range("E38").Value = "  +2.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5653"
$ws.Range("E39").Value = "  +6.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01996"
$ws.Range("E40").Value = "  +0.33%  "
$ws.Range("E41").Value = "  -0.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.973"
$ws.Range("E42").Value = "  +7.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.595"
$ws.Range("E43").Value = "  +1.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "118.36"
$ws.Range("E44").Value = "  -1.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1509"
$ws.Range("E45").Value = "  +3.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4826"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.21"
$ws.Range("E47").Value = "  +2.52%  "
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("E49").Value = "  +1.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.19"
$ws.Range("E50").Value = "  +0.80%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.58"
$ws.Range("E51").Value = "  +1.16%  "
